$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 317; existing rows 317:386 shift down to 318:387
$ws.Rows.Item(317).Insert()

# Duplicate the row that is now at 318 (the former row 317) into the new blank
# row 317 so the fixed/static columns (A,B,C,E,F,G,H,I,N,O,Q,R) and formatting
# are carried over correctly.
$ws.Range("A318:R318").Copy()
$ws.Range("A317").PasteSpecial()

# Now overwrite the new weekly record's changed values in row 317.
$ws.Range("D317").Value = 44785
$ws.Range("J317").Value = 230
$ws.Range("K317").Value = 4000
$ws.Range("L317").Value = 4500
$ws.Range("M317").Value = 4239
$ws.Range("P317").Value = 1413
